$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (MA03): weather changes from "Normal" to "Gewitter"
$ws.Range("F4").Value = "Gewitter"

# Row 4 (MA03): regenschirm changes from 0 (number) to TRUE (boolean)
$ws.Range("G4").Value = $true

# Row 5 (MA04): regenschirm changes from FALSE (boolean) to 0 (number)
$ws.Range("G5").Value = 0
